$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix i2c (issue #48)
# * lower R48 to 2k2 -> FP LEDs current halved (12*2 -> 12)
$ws.Range("F13").Formula = "=12"

# * connect IC10 directly to EEM0_I2C (before T6, T9) -> move its 12mA draw from F (P3V3) to L (P3V3MP)
$ws.Range("L13").Value = 12

# * connect IC10 and config pins to P3V3_MP -> move AT25SF081 current from F (P3V3) to L (P3V3MP)
$ws.Range("F29").ClearContents()
$ws.Range("L29").Value = 20

# restore the active selection as left by the author
$ws.Range("L32").Select()
